$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 193.6
$ws.Range("I2").Value = 191.5
$ws.Range("J2").Value = 195
$ws.Range("K2").Value = 191.5
$ws.Range("L2").Value = 195
$ws.Range("M2").Value = -78.5
$ws.Range("N2").Value = -421
# row 5
$ws.Range("H5").Value = 55.8
$ws.Range("I5").Value = 55.8
$ws.Range("K5").Value = 55.8
$ws.Range("M5").Value = 59.2
# row 18
$ws.Range("H18").Value = 5135.1665
$ws.Range("I18").Value = 2955.5293
$ws.Range("K18").Value = 2955.5293
$ws.Range("M18").Value = -2671.5293
# row 40
$ws.Range("H40").Value = 6333.3335
$ws.Range("J40").Value = 8333.333000000001
$ws.Range("L40").Value = 8333.333000000001
$ws.Range("N40").Value = -8683.333000000001
# row 112
$ws.Range("H112").Value = 2192.5
$ws.Range("J112").Value = 2198.6667
$ws.Range("L112").Value = 6596.000100000001
$ws.Range("N112").Value = -8812.000100000001
# row 116
$ws.Range("H116").Value = 7105.8696
$ws.Range("I116").Value = 5287.6665
$ws.Range("K116").Value = 5287.6665
$ws.Range("M116").Value = -1845.6665
# row 138
$ws.Range("H138").Value = 2749.4736
$ws.Range("I138").Value = 1282.7368
$ws.Range("J138").Value = 4216.2104
$ws.Range("K138").Value = 3848.2104
$ws.Range("L138").Value = 12648.6312
$ws.Range("M138").Value = 1291.7896
$ws.Range("N138").Value = -22928.6312

$ws = $wb.Worksheets.Item("ARM")
# row 35
$ws.Range("H35").Value = 2566.3333
$ws.Range("I35").Value = 2566.3333
$ws.Range("K35").Value = 2566.3333
$ws.Range("M35").Value = -2160.3333
# row 61
$ws.Range("H61").Value = 3187.8462
$ws.Range("I61").Value = 2161
$ws.Range("K61").Value = 2161
$ws.Range("M61").Value = -1949
# row 63
$ws.Range("H63").Value = 16099.8
$ws.Range("I63").Value = 20000
$ws.Range("J63").Value = 15124.75
$ws.Range("K63").Value = 20000
$ws.Range("L63").Value = 15124.75
$ws.Range("M63").Value = -19314
$ws.Range("N63").Value = -16496.75
# row 66
$ws.Range("H66").Value = 16099.8
$ws.Range("I66").Value = 20000
$ws.Range("J66").Value = 15124.75
$ws.Range("K66").Value = 100000
$ws.Range("L66").Value = 75623.75
$ws.Range("M66").Value = -96568
$ws.Range("N66").Value = -82487.75
# row 97
$ws.Range("H97").Value = 4862
$ws.Range("I97").Value = 2899.25
$ws.Range("J97").Value = 6824.75
$ws.Range("K97").Value = 2899.25
$ws.Range("L97").Value = 6824.75
$ws.Range("M97").Value = -2403.25
$ws.Range("N97").Value = -7816.75
# row 132
$ws.Range("H132").Value = 3216
$ws.Range("I132").Value = 2324.25
$ws.Range("K132").Value = 6972.75
$ws.Range("M132").Value = -4442.75
# row 136
$ws.Range("H136").Value = 3187.8462
$ws.Range("I136").Value = 2161
$ws.Range("K136").Value = 6483
$ws.Range("M136").Value = -3933
# row 139
$ws.Range("H139").Value = 60000
$ws.Range("I139").Value = 60000
$ws.Range("J139").Value = 60000
$ws.Range("K139").Value = 60000
$ws.Range("L139").Value = 60000
$ws.Range("M139").Value = -54860
$ws.Range("N139").Value = -70280

$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Range("H99").Value = 58374.75
$ws.Range("I99").Value = 35334
$ws.Range("K99").Value = 35334
$ws.Range("M99").Value = -33836
# row 107
$ws.Range("H107").Value = 7584
$ws.Range("I107").Value = 7063.077
$ws.Range("K107").Value = 7063.077
$ws.Range("M107").Value = -5143.077
# row 134
$ws.Range("H134").Value = 1448.4
$ws.Range("I134").Value = 1194.7858
$ws.Range("K134").Value = 3584.3574
$ws.Range("M134").Value = -1049.3574

$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 331.2
$ws.Range("I7").Value = 172.66667
$ws.Range("K7").Value = 172.66667
$ws.Range("M7").Value = -59.66667000000001
# row 22
$ws.Range("H22").Value = 437.08334
$ws.Range("I22").Value = 520.4
$ws.Range("J22").Value = 377.57144
$ws.Range("K22").Value = 520.4
$ws.Range("L22").Value = 377.57144
$ws.Range("M22").Value = -170.4
$ws.Range("N22").Value = -1077.57144
# row 31
$ws.Range("H31").Value = 2496.8
$ws.Range("I31").Value = 1902.25
$ws.Range("K31").Value = 1902.25
$ws.Range("M31").Value = -1607.25
# row 34
$ws.Range("H34").Value = 2496.8
$ws.Range("I34").Value = 1902.25
$ws.Range("K34").Value = 1902.25
$ws.Range("M34").Value = -1700.25
# row 62
$ws.Range("H62").Value = 71434530
$ws.Range("I62").Value = 8425
$ws.Range("K62").Value = 8425
$ws.Range("M62").Value = -7801
# row 65
$ws.Range("H65").Value = 71434530
$ws.Range("I65").Value = 8425
$ws.Range("K65").Value = 42125
$ws.Range("M65").Value = -39005
# row 107
$ws.Range("H107").Value = 7974
$ws.Range("I107").Value = 467.83334
$ws.Range("J107").Value = 13603.625
$ws.Range("K107").Value = 467.83334
$ws.Range("L107").Value = 13603.625
$ws.Range("M107").Value = 1452.16666
$ws.Range("N107").Value = -17443.625
# row 134
$ws.Range("H134").Value = 2213.65
$ws.Range("I134").Value = 1460.3125
$ws.Range("J134").Value = 5227
$ws.Range("K134").Value = 4380.9375
$ws.Range("L134").Value = 15681
$ws.Range("M134").Value = -1845.9375
$ws.Range("N134").Value = -20751

$ws = $wb.Worksheets.Item("CUL")
# row 92
$ws.Range("H92").Value = 799
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 799
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2397
$ws.Range("M92").Value = ""
$ws.Range("N92").Value = -4893
# row 95
$ws.Range("H95").Value = 4990
$ws.Range("I95").Value = 4990
$ws.Range("K95").Value = 14970
$ws.Range("M95").Value = -12911
# row 97
$ws.Range("H97").Value = 298.5
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").Value = ""
# row 101
$ws.Range("H101").Value = 6666
$ws.Range("J101").Value = 6666
$ws.Range("L101").Value = 19998
$ws.Range("N101").Value = -24866
# row 102
$ws.Range("H102").Value = 4250
$ws.Range("I102").Value = 4000
$ws.Range("J102").Value = 4500
$ws.Range("K102").Value = 12000
$ws.Range("L102").Value = 13500
$ws.Range("M102").Value = -9566
$ws.Range("N102").Value = -18368

$ws = $wb.Worksheets.Item("GSM")
# row 18
$ws.Range("H18").Value = 5000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 5000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = -5586
# row 43
$ws.Range("H43").Value = 18875
$ws.Range("I43").Value = 18875
$ws.Range("K43").Value = 18875
$ws.Range("M43").Value = -18724
# row 97
$ws.Range("H97").Value = 17316.666
$ws.Range("J97").Value = 25497.5
$ws.Range("L97").Value = 25497.5
$ws.Range("N97").Value = -26489.5
# row 126
$ws.Range("H126").Value = 1402.1666
$ws.Range("I126").Value = 786
$ws.Range("K126").Value = 2358
$ws.Range("M126").Value = 112

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 18182716
$ws.Range("J22").Value = 1993
$ws.Range("L22").Value = 1993
$ws.Range("N22").Value = -2583
# row 27
$ws.Range("H27").Value = 18182716
$ws.Range("J27").Value = 1993
$ws.Range("L27").Value = 1993
$ws.Range("N27").Value = -2207
# row 46
$ws.Range("H46").Value = 3174.875
$ws.Range("J46").Value = 3997.0908
$ws.Range("L46").Value = 3997.0908
$ws.Range("N46").Value = -4373.0908
# row 93
$ws.Range("H93").Value = 1108.2222
$ws.Range("I93").Value = 601
$ws.Range("J93").Value = 1742.25
$ws.Range("K93").Value = 601
$ws.Range("L93").Value = 1742.25
$ws.Range("M93").Value = 647
$ws.Range("N93").Value = -4238.25
# row 100
$ws.Range("H100").Value = 185764.83
$ws.Range("I100").Value = 185764.83
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 185764.83
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -185223.83
$ws.Range("N100").Value = ""

$ws = $wb.Worksheets.Item("WVR")
# row 100
$ws.Range("H100").Value = 577.6
$ws.Range("I100").Value = 522
$ws.Range("K100").Value = 1044
$ws.Range("M100").Value = -503
# row 107
$ws.Range("H107").Value = 987
$ws.Range("J107").Value = 999
$ws.Range("L107").Value = 2997
$ws.Range("N107").Value = -6837
# row 113
$ws.Range("H113").Value = 1095.7142
$ws.Range("I113").Value = 482.55554
$ws.Range("J113").Value = 2199.4
$ws.Range("K113").Value = 1447.66662
$ws.Range("L113").Value = 6598.200000000001
$ws.Range("M113").Value = 722.33338
$ws.Range("N113").Value = -10938.2
# row 132
$ws.Range("H132").Value = 2876.9167
$ws.Range("I132").Value = 2876.9167
$ws.Range("K132").Value = 8630.750100000001
$ws.Range("M132").Value = -6100.750100000001
